$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualización desde MV -datos-": correct the previous last row's D value
# (128 -> 127.9) and append the new monthly observation (01-09-2021).

$ws.Range("D69").Value = 127.9

# A70 must hold the literal text "01-09-2021" (a shared string), not an
# auto-converted date serial, so force Text formatting before assigning the
# value, then drop back to the sheet's normal (unformatted) style so no
# stray number format lingers on the cell.
$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = "01-09-2021"
$ws.Range("A70").Style = "Normal"

$ws.Range("B70").Value = 127.7
$ws.Range("C70").Value = 96.5
$ws.Range("D70").Value = 128.8
